$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Shorten the Observational data definition (row 36, column C)
$ws.Range("C36").Value = "Data collected from a study where researchers are observing the effect of an intervention without manipulating who is exposed to the intervention."

# 2. Fix typo "identifable" -> "identifiable" in the Confidentiality definition (row 11, column C)
$ws.Range("C11").Value = "Confidentiality concerns data, ensuring participants agree to how their private and identifiable information will be managed and disseminated."

# 3. Clear the custom fill formatting on rows 3, 42, 43, 51 (revert to default style)
$ws.Rows(3).ClearFormats()
$ws.Rows(42).ClearFormats()
$ws.Rows(43).ClearFormats()
$ws.Rows(51).ClearFormats()

# 4. Update the active cell selection
$ws.Range("E11").Select()
